$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6263.0835
$ws.Range("I28").Value = 779.1429000000001
$ws.Range("J28").Value = 13940.6
$ws.Range("K28").Value = 779.1429000000001
$ws.Range("L28").Value = 13940.6
$ws.Range("M28").Value = -294.1429000000001
$ws.Range("N28").Value = -14910.6

$ws.Range("H39").Value = 560.9167
$ws.Range("J39").Value = 2925
$ws.Range("L39").Value = 8775
$ws.Range("N39").Value = -9367

$ws.Range("H53").Value = 340.4
$ws.Range("J53").Value = 812.2
$ws.Range("L53").Value = 812.2
$ws.Range("N53").Value = -2086.2

$ws.Range("H98").Value = 1912.8182
$ws.Range("I98").Value = 1912.8182
$ws.Range("K98").Value = 1912.8182
$ws.Range("M98").Value = -414.8181999999999

$ws.Range("H107").Value = 2285.2
$ws.Range("I107").Value = 1124.25
$ws.Range("J107").Value = 4026.625
$ws.Range("K107").Value = 1124.25
$ws.Range("L107").Value = 4026.625
$ws.Range("M107").Value = 795.75
$ws.Range("N107").Value = -7866.625

$ws.Range("H116").Value = 18211.25
$ws.Range("I116").Value = 23778.4
$ws.Range("K116").Value = 23778.4
$ws.Range("M116").Value = -20336.4

$ws.Range("H122").Value = 1912.8182
$ws.Range("I122").Value = 1912.8182
$ws.Range("K122").Value = 5738.4546
$ws.Range("M122").Value = -3288.4546

$ws.Range("H127").Value = 7816.1055
$ws.Range("I127").Value = 7055.625
$ws.Range("J127").Value = 11872
$ws.Range("K127").Value = 21166.875
$ws.Range("L127").Value = 35616
$ws.Range("M127").Value = -16206.875
$ws.Range("N127").Value = -45536

$ws.Range("H132").Value = 6522.115
$ws.Range("I132").Value = 5122.091
$ws.Range("J132").Value = 14222.25
$ws.Range("K132").Value = 15366.273
$ws.Range("L132").Value = 42666.75
$ws.Range("M132").Value = -12836.273
$ws.Range("N132").Value = -47726.75

$ws.Range("H133").Value = 79000
$ws.Range("J133").Value = 79000
$ws.Range("L133").Value = 79000
$ws.Range("N133").Value = -89120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7955.5
$ws.Range("I61").Value = 16499.5
$ws.Range("J61").Value = 5107.5
$ws.Range("K61").Value = 16499.5
$ws.Range("L61").Value = 5107.5
$ws.Range("M61").Value = -16287.5
$ws.Range("N61").Value = -5531.5

$ws.Range("H63").Value = 1919.0769
$ws.Range("I63").Value = 1374.9
$ws.Range("J63").Value = 3733
$ws.Range("K63").Value = 1374.9
$ws.Range("L63").Value = 3733
$ws.Range("M63").Value = -688.9000000000001
$ws.Range("N63").Value = -5105

$ws.Range("H66").Value = 1919.0769
$ws.Range("I66").Value = 1374.9
$ws.Range("J66").Value = 3733
$ws.Range("K66").Value = 6874.5
$ws.Range("L66").Value = 18665
$ws.Range("M66").Value = -3442.5
$ws.Range("N66").Value = -25529

$ws.Range("H122").Value = 6835.25
$ws.Range("I122").Value = 6695
$ws.Range("K122").Value = 20085
$ws.Range("M122").Value = -17635

$ws.Range("H132").Value = 42532.957
$ws.Range("I132").Value = 2801.0833
$ws.Range("J132").Value = 82264.836
$ws.Range("K132").Value = 8403.249899999999
$ws.Range("L132").Value = 246794.508
$ws.Range("M132").Value = -5873.249899999999
$ws.Range("N132").Value = -251854.508

$ws.Range("H136").Value = 7955.5
$ws.Range("I136").Value = 16499.5
$ws.Range("J136").Value = 5107.5
$ws.Range("K136").Value = 49498.5
$ws.Range("L136").Value = 15322.5
$ws.Range("M136").Value = -46948.5
$ws.Range("N136").Value = -20422.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16360.454
$ws.Range("I58").Value = 10831.333
$ws.Range("K58").Value = 10831.333
$ws.Range("M58").Value = -10628.333

$ws.Range("H132").Value = 7254.222
$ws.Range("I132").Value = 6607.4116
$ws.Range("K132").Value = 19822.2348
$ws.Range("M132").Value = -17292.2348

$ws.Range("H134").Value = 472797.2
$ws.Range("I134").Value = 1164952.9
$ws.Range("J134").Value = 22895.95
$ws.Range("K134").Value = 3494858.7
$ws.Range("L134").Value = 68687.85000000001
$ws.Range("M134").Value = -3492323.7
$ws.Range("N134").Value = -73757.85000000001

$ws.Range("H136").Value = 16360.454
$ws.Range("I136").Value = 10831.333
$ws.Range("K136").Value = 32493.999
$ws.Range("M136").Value = -29943.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 13999.833
$ws.Range("J62").Value = 18000
$ws.Range("L62").Value = 54000
$ws.Range("N62").Value = -55372

$ws.Range("H65").Value = 13999.833
$ws.Range("J65").Value = 18000
$ws.Range("L65").Value = 162000
$ws.Range("N65").Value = -168864

$ws.Range("H104").Value = 1458.25
$ws.Range("I104").Value = 499
$ws.Range("K104").Value = 1497
$ws.Range("M104").Value = 1124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11924
$ws.Range("I102").Value = 9315.833000000001
$ws.Range("J102").Value = 19748.5
$ws.Range("K102").Value = 9315.833000000001
$ws.Range("L102").Value = 19748.5
$ws.Range("M102").Value = -7693.833000000001
$ws.Range("N102").Value = -22992.5

$ws.Range("H132").Value = 3811.3635
$ws.Range("I132").Value = 3946.1482
$ws.Range("J132").Value = 3204.8333
$ws.Range("K132").Value = 11838.4446
$ws.Range("L132").Value = 9614.499899999999
$ws.Range("M132").Value = -9308.444600000001
$ws.Range("N132").Value = -14674.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6437.9546
$ws.Range("I61").Value = 4859.2
$ws.Range("K61").Value = 4859.2
$ws.Range("M61").Value = -4657.2

$ws.Range("H82").Value = 2948.913
$ws.Range("I82").Value = 1218.2
$ws.Range("K82").Value = 1218.2
$ws.Range("M82").Value = -857.2

$ws.Range("H85").Value = 2948.913
$ws.Range("I85").Value = 1218.2
$ws.Range("K85").Value = 1218.2
$ws.Range("M85").Value = 29.79999999999995

$ws.Range("H113").Value = 6437.9546
$ws.Range("I113").Value = 4859.2
$ws.Range("K113").Value = 4859.2
$ws.Range("M113").Value = -2689.2

$ws.Range("H121").Value = 58868.332
$ws.Range("J121").Value = 58868.332
$ws.Range("L121").Value = 58868.332
$ws.Range("N121").Value = -62362.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7510.2666
$ws.Range("I62").Value = 5902.636
$ws.Range("J62").Value = 11931.25
$ws.Range("K62").Value = 5902.636
$ws.Range("L62").Value = 11931.25
$ws.Range("M62").Value = -5278.636
$ws.Range("N62").Value = -13179.25

$ws.Range("H65").Value = 7510.2666
$ws.Range("I65").Value = 5902.636
$ws.Range("J65").Value = 11931.25
$ws.Range("K65").Value = 29513.18
$ws.Range("L65").Value = 59656.25
$ws.Range("M65").Value = -26393.18
$ws.Range("N65").Value = -65896.25

$ws.Range("H81").Value = 1451.4375
$ws.Range("I81").Value = 1079
$ws.Range("J81").Value = 1823.875
$ws.Range("K81").Value = 2158
$ws.Range("L81").Value = 3647.75
$ws.Range("M81").Value = -1097
$ws.Range("N81").Value = -5769.75

$ws.Range("H84").Value = 1451.4375
$ws.Range("I84").Value = 1079
$ws.Range("J84").Value = 1823.875
$ws.Range("K84").Value = 10790
$ws.Range("L84").Value = 18238.75
$ws.Range("M84").Value = -5486
$ws.Range("N84").Value = -28846.75

$ws.Range("H96").Value = 1571.9231
$ws.Range("I96").Value = 1414.9166
$ws.Range("J96").Value = 1706.5
$ws.Range("K96").Value = 1414.9166
$ws.Range("L96").Value = 1706.5
$ws.Range("M96").Value = -41.91660000000002
$ws.Range("N96").Value = -4452.5

$ws.Range("H121").Value = 10000
$ws.Range("J121").Value = 10000
$ws.Range("L121").Value = 10000
$ws.Range("N121").Value = -13494

$ws.Range("H132").Value = 11114.071
$ws.Range("J132").Value = 22854.857
$ws.Range("L132").Value = 68564.571
$ws.Range("N132").Value = -73624.571

$ws.Range("H136").Value = 38807.668
$ws.Range("I136").Value = 43146
$ws.Range("K136").Value = 129438
$ws.Range("M136").Value = -126888
